$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from row 374 (A column) down through the newly added rows
$ws.Range("A374").Copy()
$ws.Range("A375:A385").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A375").Value = 44449
$ws.Range("B375").Value = 1
$ws.Range("C375").Value = 1
$ws.Range("D375").Value = 18.93939393939394

$ws.Range("A376").Value = 44450
$ws.Range("B376").Value = 0
$ws.Range("C376").Value = 1
$ws.Range("D376").Value = 18.93939393939394

$ws.Range("A377").Value = 44451
$ws.Range("B377").Value = 1
$ws.Range("C377").Value = 2
$ws.Range("D377").Value = 37.87878787878788

$ws.Range("A378").Value = 44452
$ws.Range("B378").Value = 3
$ws.Range("C378").Value = 5
$ws.Range("D378").Value = 94.6969696969697

$ws.Range("A379").Value = 44453
$ws.Range("B379").Value = 0
$ws.Range("C379").Value = 5
$ws.Range("D379").Value = 94.6969696969697

$ws.Range("A380").Value = 44454
$ws.Range("B380").Value = 0
$ws.Range("C380").Value = 5
$ws.Range("D380").Value = 94.6969696969697

$ws.Range("A381").Value = 44455
$ws.Range("B381").Value = 0
$ws.Range("C381").Value = 5
$ws.Range("D381").Value = 94.6969696969697

$ws.Range("A382").Value = 44456
$ws.Range("B382").Value = 1
$ws.Range("C382").Value = 5
$ws.Range("D382").Value = 94.6969696969697

$ws.Range("A383").Value = 44457
$ws.Range("B383").Value = 1
$ws.Range("C383").Value = 6
$ws.Range("D383").Value = 113.6363636363636

$ws.Range("A384").Value = 44458
$ws.Range("B384").Value = 0
$ws.Range("C384").Value = 5
$ws.Range("D384").Value = 94.6969696969697

$ws.Range("A385").Value = 44459
$ws.Range("B385").Value = 0
$ws.Range("C385").Value = 2
$ws.Range("D385").Value = 37.87878787878788

